# Applies the edit described by the commit diff:
#  1. In the intro paragraph, add spell-check proofing marks (w:proofErr
#     spellStart/spellEnd) around "gyors", "internetkapcsolatot", and move
#     the stray "_GoBack" bookmark away from in front of "gyors".
#  2. In the "6. Biztonsag" section's body paragraph, add proofing marks
#     around "halozat", "tuzfallal", "es", "titkositassal", drop the
#     trailing VPN sentence, and re-anchor the "_GoBack" bookmark at the
#     end of the paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part A: "6. Biztonsag" paragraph -- do this FIRST so the freshly
# inserted bookmarkStart/bookmarkEnd picks up id 0 (the only "_GoBack"
# bookmark still in the document at this point is the original one that
# part B below is about to remove).
# ---------------------------------------------------------------------
$rFindB = $d.Content.Duplicate
$foundB = $rFindB.Find.Execute("A hálózat minden szegmense", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundB) {
    throw "Could not locate the 'Biztonsag' paragraph"
}
$paraB = $rFindB.Paragraphs(1).Range
$rB = $d.Range($paraB.Start, $paraB.End - 1)

$xmlB = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hálózat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tűzfallal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>és</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>titkosítással</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> van védve, hogy megakadályozza az illetéktelen hozzáférést és biztosítsa az adatok védelmét.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rB.InsertXML($xmlB)

# ---------------------------------------------------------------------
# Part B: intro paragraph ("Az iskolai halozat kiepitese gyors es
# megbizhato internetkapcsolatot ...").
# ---------------------------------------------------------------------
$rFindA = $d.Content.Duplicate
$foundA = $rFindA.Find.Execute("Az iskolai hálózat kiépítése", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundA) {
    throw "Could not locate the intro paragraph"
}
$paraA = $rFindA.Paragraphs(1).Range
$rA = $d.Range($paraA.Start, $paraA.End - 1)

$xmlA = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Az </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>iskolai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hálózat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>kiépítése</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gyors</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>és</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>megbízható</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>internetkapcsolatot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> biztosítja, hanem lehetőséget ad a tanároknak, diákoknak és vendégeknek, hogy hatékonyan használják a különböző digitális eszközöket. A hálózat </w:t></w:r><w:r><w:t xml:space="preserve">fa topológiára épül, és ethernet kábelekkel biztosítják a stabil vezetékes kapcsolatot. A vezeték nélküli hozzáférést Ubiquiti UniFi U6+ access pointokkal valósítják meg, amelyek a legújabb Wi-Fi 6 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>szabványt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>t</w:t></w:r><w:r><w:t>á</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>o</w:t></w:r><w:r><w:t>g</w:t></w:r><w:r><w:t>a</w:t></w:r><w:r><w:t>t</w:t></w:r><w:r><w:t>j</w:t></w:r><w:r><w:t>á</w:t></w:r><w:r><w:t>k</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rA.InsertXML($xmlA)
